# Fruta / hortaliza, semanal
#
# Weekly price-sheet update: a new daily/weekly record for
# "Feria Lagunitas de Puerto Montt - Ajo" (Chino, Primera) is inserted
# ahead of the existing historical rows, pushing the prior rows
# (203..326) down by one (204..327) and extending the used range from
# A1:R326 to A1:R327.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 203, shifting everything below
# (203..326) down to (204..327).
$ws.Rows(203).Insert()

# Populate the newly inserted row with the new price record. All the
# "constant" columns (Mercado ID/Mercado/Region/Codreg/Categoria
# ID/Categoria/Unidad/Origen/Kg o Unidades/Clasificacion) match every
# other row in this sheet.
$ws.Range("A203").Value = 4
$ws.Range("B203").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C203").Value = "Los Lagos"
$ws.Range("D203").Value = 44824
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = 100112003
$ws.Range("G203").Value = "Ajo"
$ws.Range("H203").Value = "Chino"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 30
$ws.Range("K203").Value = 23000
$ws.Range("L203").Value = 23000
$ws.Range("M203").Value = 23000
$ws.Range("N203").Value = "$/caja 10 kilos"
$ws.Range("O203").Value = "China"
$ws.Range("P203").Value = 2300
$ws.Range("Q203").Value = 10
$ws.Range("R203").Value = "Hortaliza"
